$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-23 from 45180 to 45181
$ws.Range("C2:C23").Value = 45181
